# Fix some flaws of POR calculation and CE member selection.
#
# For each 6-row block (rows 1-6, 7-12, 13-18, ... 55-60) in columns A:C:
#   - the first two rows of the block get brand new figures (newly selected
#     CE members), and
#   - the POR-calculation rows that used to sit in positions 2 and 3 of the
#     block shift down to positions 3 and 4 (the old position-4 figures, an
#     earlier calculation error, are discarded).
# Rows 5 and 6 of every block (and the trailing all-zero rows 61-72) are
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the first row of each block (columns A, B, C).
$newTop = @{
    1  = @(9235, 7162, 7077)
    7  = @(8467, 6405, 6598)
    13 = @(9386, 6289, 6385)
    19 = @(7933, 4609, 4624)
    25 = @(10566, 7921, 8019)
    31 = @(14563, 9290, 9047)
    37 = @(22846, 13406, 11278)
    43 = @(13791, 13088, 11657)
    49 = @(9795, 9224, 8784)
    55 = @(9961, 8720, 8305)
}

# New values for the second row of each block (columns A, B, C).
$newSecond = @{
    1  = @(6668, 2856, 2973)
    7  = @(5513, 2471, 2556)
    13 = @(5370, 2316, 2309)
    19 = @(4364, 1440, 1579)
    25 = @(5873, 2606, 2656)
    31 = @(6580, 3214, 3375)
    37 = @(8977, 3985, 4141)
    43 = @(8247, 4344, 4509)
    49 = @(7327, 3234, 3381)
    55 = @(5831, 2479, 2583)
}

$blockStarts = @(1, 7, 13, 19, 25, 31, 37, 43, 49, 55)

foreach ($r0 in $blockStarts) {
    # Capture the current (pre-edit) values of rows r0+1 and r0+2 before
    # anything in this block gets overwritten. NOTE: the COM bridge's
    # parameterized `Value` getter must be invoked with explicit `()` or it
    # resolves to the property descriptor instead of the cell's value.
    $oldRow2 = @($ws.Cells.Item($r0 + 1, 1).Value(), $ws.Cells.Item($r0 + 1, 2).Value(), $ws.Cells.Item($r0 + 1, 3).Value())
    $oldRow3 = @($ws.Cells.Item($r0 + 2, 1).Value(), $ws.Cells.Item($r0 + 2, 2).Value(), $ws.Cells.Item($r0 + 2, 3).Value())

    # Row r0+3 <- old row r0+2 ; row r0+2 <- old row r0+1 (shift down by one).
    for ($c = 1; $c -le 3; $c++) {
        $ws.Cells.Item($r0 + 3, $c).Value = $oldRow3[$c - 1]
        $ws.Cells.Item($r0 + 2, $c).Value = $oldRow2[$c - 1]
    }

    # Row r0 and r0+1 <- brand new figures.
    $top = $newTop[$r0]
    $second = $newSecond[$r0]
    for ($c = 1; $c -le 3; $c++) {
        $ws.Cells.Item($r0, $c).Value = $top[$c - 1]
        $ws.Cells.Item($r0 + 1, $c).Value = $second[$c - 1]
    }
}

Write-Output "done"
